# Append a new data row (row 3) to Sheet1 of the users_engagement workbook,
# mirroring the row-2 layout: user_id, username, level, last_message_date,
# last_response, response_status, level_3_ai_response, subscription_checked,
# level_4_reminder_sent, decision, notes, first_added_date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 3

$ws.Cells.Item($row, 1).Value  = 1594575966
$ws.Cells.Item($row, 2).Value  = "user_1594575966"
$ws.Cells.Item($row, 3).Value  = 0
$ws.Cells.Item($row, 4).Value  = "2026-02-14T04:59:35.710181+00:00"
$ws.Cells.Item($row, 5).Value  = ""
$ws.Cells.Item($row, 6).Value  = ""
$ws.Cells.Item($row, 7).Value  = ""
$ws.Cells.Item($row, 8).Value  = $false
$ws.Cells.Item($row, 9).Value  = $false
$ws.Cells.Item($row, 10).Value = ""
$ws.Cells.Item($row, 11).Value = "Added during extraction"

# Column L holds a plain ISO date string ("2026-02-14"), not a real date
# serial. Enter it with a leading apostrophe so Excel stores it as text
# instead of auto-converting it to a date value, then reset the style back
# to Normal so no date number-format sticks to the cell.
$ws.Cells.Item($row, 12).Value = "'2026-02-14"
$ws.Cells.Item($row, 12).Style = "Normal"
